$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (shifts existing update_usr_id/update_time columns
# from J/K to K/L) to make room for the new tenant_id column.
$ws.Range("J1:J2").EntireColumn.Insert()

# Populate the new tenant_id column's header (comment) and detail (model) cells.
$ws.Range("J1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
$ws.Range("J2").Value = '<%=model.tenant_id_lbl%>'
